$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.907.79"
$ws.Range("E2").Value = "  +2.44%  "

# Row 3
$ws.Range("D3").Value = "2.501.35"
$ws.Range("E3").Value = "  +0.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.90"
$ws.Range("E5").Value = "  +3.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.44"
$ws.Range("E6").Value = "  +9.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +1.27%  "

# Row 9
$ws.Range("D9").Value = "2.516.28"
$ws.Range("E9").Value = "  +0.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.79"
$ws.Range("E10").Value = "  +6.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0992"
$ws.Range("E11").Value = "  +0.94%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  +2.90%  "

# Row 13
$ws.Range("E13").Value = "  +1.33%  "

# Row 14
$ws.Range("D14").Value = "2.937.37"
$ws.Range("E14").Value = "  +0.19%  "

# Row 15
$ws.Range("D15").Value = "56.969.31"
$ws.Range("E15").Value = "  +2.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.39"
$ws.Range("E16").Value = "  +3.93%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("D18").Value = "2.512.58"
$ws.Range("E18").Value = "  +0.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  +4.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.35"
$ws.Range("E20").Value = "  +3.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.49"
$ws.Range("E21").Value = "  +1.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("E23").Value = "  +4.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.08"
$ws.Range("E24").Value = "  +2.21%  "

# Row 25
$ws.Range("E25").Value = "  +1.39%  "

# Row 26
$ws.Range("E26").Value = "  -0.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.84%  "

# Row 28
$ws.Range("D28").Value = "2.610.58"
$ws.Range("E28").Value = "  +0.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"

# Row 30
$ws.Range("D30").Value = "0.0₃0819"
$ws.Range("E30").Value = "  +4.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.23%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.50"
$ws.Range("E32").Value = "  +1.99%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.45"
$ws.Range("E33").Value = "  +1.59%  "

# Row 34
$ws.Range("E34").Value = "  +3.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  +1.41%  "

# Row 36
$ws.Range("E36").Value = "  +5.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +2.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("E38").Value = "  +4.15%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.35"
$ws.Range("E39").Value = "  +0.26%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.39"
$ws.Range("E40").Value = "  +5.88%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0567"
$ws.Range("E41").Value = "  +3.06%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.52"
$ws.Range("E42").Value = "  +4.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.616"
$ws.Range("E43").Value = "  +0.84%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.94"
$ws.Range("E45").Value = "  +9.79%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "268.45"
$ws.Range("E46").Value = "  +8.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0927"
$ws.Range("E47").Value = "  +2.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  +3.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.23"
$ws.Range("E49").Value = "  +0.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.97"
$ws.Range("E50").Value = "  +2.93%  "

# Row 51
$ws.Range("D51").Value = "1.911.54"
$ws.Range("E51").Value = "  -3.08%  "
